$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44858
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = 7500
$ws.Range("N2").Value = '$/saco 25 kilos'
$ws.Range("O2").Value = 'Provincia del Elquí'
$ws.Range("P2").Value = 300

$ws.Range("D3").Value = 44377
$ws.Range("J3").Value = 520
$ws.Range("K3").Value = 12500
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 12750
$ws.Range("N3").Value = '$/saco 25 kilos'
$ws.Range("O3").Value = 'Provincia del Elquí'
$ws.Range("P3").Value = 510

$ws.Range("D4").Value = 44846
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 7500
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 7750
$ws.Range("N4").Value = '$/saco 25 kilos'
$ws.Range("O4").Value = 'Provincia del Elquí'
$ws.Range("P4").Value = 310

$ws.Range("D5").Value = 44809
$ws.Range("J5").Value = 520
$ws.Range("K5").Value = 9500
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 9750
$ws.Range("N5").Value = '$/saco 25 kilos'
$ws.Range("O5").Value = 'Provincia del Elquí'
$ws.Range("P5").Value = 390

$ws.Range("D6").Value = 44837
$ws.Range("J6").Value = 520
$ws.Range("K6").Value = 8000
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = 8500
$ws.Range("N6").Value = '$/saco 25 kilos'
$ws.Range("O6").Value = 'Provincia del Elquí'
$ws.Range("P6").Value = 340

$ws.Range("D7").Value = 44824
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 8000
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = 8500
$ws.Range("N7").Value = '$/saco 25 kilos'
$ws.Range("O7").Value = 'Provincia del Elquí'
$ws.Range("P7").Value = 340

$ws.Range("D8").Value = 44799
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 10000
$ws.Range("L8").Value = 11000
$ws.Range("M8").Value = 10500
$ws.Range("N8").Value = '$/saco 25 kilos'
$ws.Range("O8").Value = 'Provincia del Elquí'
$ws.Range("P8").Value = 420

$ws.Range("D9").Value = 44690
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 17000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 17500
$ws.Range("N9").Value = '$/saco 25 kilos'
$ws.Range("O9").Value = 'Provincia del Elquí'
$ws.Range("P9").Value = 700

$ws.Range("D10").Value = 44372
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 13000
$ws.Range("L10").Value = 14000
$ws.Range("M10").Value = 13500
$ws.Range("N10").Value = '$/saco 25 kilos'
$ws.Range("O10").Value = 'Provincia del Elquí'
$ws.Range("P10").Value = 540

$ws.Range("D11").Value = 44370
$ws.Range("J11").Value = 520
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 13500
$ws.Range("N11").Value = '$/saco 25 kilos'
$ws.Range("O11").Value = 'Provincia del Elquí'
$ws.Range("P11").Value = 540

$ws.Range("D12").Value = 44694
$ws.Range("J12").Value = 480
$ws.Range("K12").Value = 17500
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 17750
$ws.Range("N12").Value = '$/saco 25 kilos'
$ws.Range("O12").Value = 'Provincia del Elquí'
$ws.Range("P12").Value = 710

$ws.Range("D13").Value = 44881
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 8000
$ws.Range("M13").Value = 7500
$ws.Range("N13").Value = '$/saco 25 kilos'
$ws.Range("O13").Value = 'Provincia del Elquí'
$ws.Range("P13").Value = 300

$ws.Range("D14").Value = 44817
$ws.Range("J14").Value = 440
$ws.Range("K14").Value = 9000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 9500
$ws.Range("N14").Value = '$/saco 25 kilos'
$ws.Range("O14").Value = 'Provincia del Elquí'
$ws.Range("P14").Value = 380

$ws.Range("D15").Value = 44384
$ws.Range("J15").Value = 560
$ws.Range("K15").Value = 11500
$ws.Range("L15").Value = 12000
$ws.Range("M15").Value = 11750
$ws.Range("N15").Value = '$/saco 25 kilos'
$ws.Range("O15").Value = 'Provincia del Elquí'
$ws.Range("P15").Value = 470

$ws.Range("D16").Value = 44466
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 9500
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 9750
$ws.Range("N16").Value = '$/saco 25 kilos'
$ws.Range("O16").Value = 'Provincia del Elquí'
$ws.Range("P16").Value = 390

$ws.Range("D17").Value = 44883
$ws.Range("J17").Value = 380
$ws.Range("K17").Value = 7000
$ws.Range("L17").Value = 8000
$ws.Range("M17").Value = 7500
$ws.Range("N17").Value = '$/saco 25 kilos'
$ws.Range("O17").Value = 'Provincia del Elquí'
$ws.Range("P17").Value = 300

$ws.Range("D18").Value = 44781
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 11000
$ws.Range("M18").Value = 10500
$ws.Range("N18").Value = '$/saco 25 kilos'
$ws.Range("O18").Value = 'Provincia del Elquí'
$ws.Range("P18").Value = 420

$ws.Range("D19").Value = 44484
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 9000
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = 9500
$ws.Range("N19").Value = '$/saco 25 kilos'
$ws.Range("O19").Value = 'Provincia del Elquí'
$ws.Range("P19").Value = 380

$ws.Range("D20").Value = 44848
$ws.Range("J20").Value = 800
$ws.Range("K20").Value = 7000
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 7500
$ws.Range("N20").Value = '$/saco 25 kilos'
$ws.Range("O20").Value = 'Provincia del Elquí'
$ws.Range("P20").Value = 300

$ws.Range("D21").Value = 44316
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 16000
$ws.Range("L21").Value = 17000
$ws.Range("M21").Value = 16500
$ws.Range("N21").Value = '$/saco 25 kilos'
$ws.Range("O21").Value = 'Provincia del Elquí'
$ws.Range("P21").Value = 660

$ws.Range("D22").Value = 44446
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 11000
$ws.Range("L22").Value = 12000
$ws.Range("M22").Value = 11500
$ws.Range("N22").Value = '$/saco 25 kilos'
$ws.Range("O22").Value = 'Provincia del Elquí'
$ws.Range("P22").Value = 460

$ws.Range("D23").Value = 44803
$ws.Range("J23").Value = 600
$ws.Range("K23").Value = 9500
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = 9750
$ws.Range("N23").Value = '$/saco 25 kilos'
$ws.Range("O23").Value = 'Provincia del Elquí'
$ws.Range("P23").Value = 390

$ws.Range("D24").Value = 44816
$ws.Range("J24").Value = 600
$ws.Range("K24").Value = 9500
$ws.Range("L24").Value = 10000
$ws.Range("M24").Value = 9750
$ws.Range("N24").Value = '$/saco 25 kilos'
$ws.Range("O24").Value = 'Provincia del Elquí'
$ws.Range("P24").Value = 390

$ws.Range("D25").Value = 44811
$ws.Range("J25").Value = 400
$ws.Range("K25").Value = 10000
$ws.Range("L25").Value = 10500
$ws.Range("M25").Value = 10250
$ws.Range("N25").Value = '$/saco 25 kilos'
$ws.Range("O25").Value = 'Provincia del Elquí'
$ws.Range("P25").Value = 410

$ws.Range("D26").Value = 44880
$ws.Range("J26").Value = 560
$ws.Range("K26").Value = 7000
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = 7500
$ws.Range("N26").Value = '$/saco 25 kilos'
$ws.Range("O26").Value = 'Provincia del Elquí'
$ws.Range("P26").Value = 300

$ws.Range("D27").Value = 44847
$ws.Range("J27").Value = 520
$ws.Range("K27").Value = 7000
$ws.Range("L27").Value = 8000
$ws.Range("M27").Value = 7500
$ws.Range("N27").Value = '$/saco 25 kilos'
$ws.Range("O27").Value = 'Provincia del Elquí'
$ws.Range("P27").Value = 300

$ws.Range("D28").Value = 44714
$ws.Range("J28").Value = 400
$ws.Range("K28").Value = 14000
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = 14500
$ws.Range("N28").Value = '$/saco 25 kilos'
$ws.Range("O28").Value = 'Provincia de Limarí'
$ws.Range("P28").Value = 580

$ws.Range("D29").Value = 44756
$ws.Range("J29").Value = 400
$ws.Range("K29").Value = 14000
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = 14500
$ws.Range("N29").Value = '$/saco 25 kilos'
$ws.Range("O29").Value = 'Provincia del Elquí'
$ws.Range("P29").Value = 580

$ws.Range("D30").Value = 44825
$ws.Range("J30").Value = 440
$ws.Range("K30").Value = 8000
$ws.Range("L30").Value = 9000
$ws.Range("M30").Value = 8500
$ws.Range("N30").Value = '$/saco 25 kilos'
$ws.Range("O30").Value = 'Provincia del Elquí'
$ws.Range("P30").Value = 340

$ws.Range("D31").Value = 44376
$ws.Range("J31").Value = 400
$ws.Range("K31").Value = 12000
$ws.Range("L31").Value = 13000
$ws.Range("M31").Value = 12500
$ws.Range("N31").Value = '$/saco 25 kilos'
$ws.Range("O31").Value = 'Provincia del Elquí'
$ws.Range("P31").Value = 500

$ws.Range("D32").Value = 44855
$ws.Range("J32").Value = 540
$ws.Range("K32").Value = 7000
$ws.Range("L32").Value = 8000
$ws.Range("M32").Value = 7500
$ws.Range("N32").Value = '$/saco 25 kilos'
$ws.Range("O32").Value = 'Provincia del Elquí'
$ws.Range("P32").Value = 300

$ws.Range("D33").Value = 44798
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 10500
$ws.Range("L33").Value = 11000
$ws.Range("M33").Value = 10750
$ws.Range("N33").Value = '$/saco 25 kilos'
$ws.Range("O33").Value = 'Provincia del Elquí'
$ws.Range("P33").Value = 430

$ws.Range("D34").Value = 44356
$ws.Range("J34").Value = 500
$ws.Range("K34").Value = 13000
$ws.Range("L34").Value = 14000
$ws.Range("M34").Value = 13500
$ws.Range("N34").Value = '$/saco 25 kilos'
$ws.Range("O34").Value = 'Provincia de Limarí'
$ws.Range("P34").Value = 540

$ws.Range("D35").Value = 44721
$ws.Range("J35").Value = 500
$ws.Range("K35").Value = 14500
$ws.Range("L35").Value = 15000
$ws.Range("M35").Value = 14750
$ws.Range("N35").Value = '$/saco 25 kilos'
$ws.Range("O35").Value = 'Provincia de Limarí'
$ws.Range("P35").Value = 590

$ws.Range("D36").Value = 44425
$ws.Range("J36").Value = 400
$ws.Range("K36").Value = 11500
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = 11750
$ws.Range("N36").Value = '$/saco 25 kilos'
$ws.Range("O36").Value = 'Provincia del Elquí'
$ws.Range("P36").Value = 470

$ws.Range("D37").Value = 44827
$ws.Range("J37").Value = 700
$ws.Range("K37").Value = 8000
$ws.Range("L37").Value = 9000
$ws.Range("M37").Value = 8500
$ws.Range("N37").Value = '$/saco 25 kilos'
$ws.Range("O37").Value = 'Provincia del Elquí'
$ws.Range("P37").Value = 340

$ws.Range("D38").Value = 44873
$ws.Range("J38").Value = 540
$ws.Range("K38").Value = 6000
$ws.Range("L38").Value = 7000
$ws.Range("M38").Value = 6500
$ws.Range("N38").Value = '$/saco 25 kilos'
$ws.Range("O38").Value = 'Provincia del Elquí'
$ws.Range("P38").Value = 260

$ws.Range("D39").Value = 44797
$ws.Range("J39").Value = 1000
$ws.Range("K39").Value = 11000
$ws.Range("L39").Value = 12000
$ws.Range("M39").Value = 11500
$ws.Range("N39").Value = '$/saco 25 kilos'
$ws.Range("O39").Value = 'Provincia del Elquí'
$ws.Range("P39").Value = 460

$ws.Range("D40").Value = 44876
$ws.Range("J40").Value = 460
$ws.Range("K40").Value = 6000
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = 6500
$ws.Range("N40").Value = '$/saco 25 kilos'
$ws.Range("O40").Value = 'Provincia del Elquí'
$ws.Range("P40").Value = 260

$ws.Range("D41").Value = 44386
$ws.Range("J41").Value = 500
$ws.Range("K41").Value = 11000
$ws.Range("L41").Value = 12000
$ws.Range("M41").Value = 11500
$ws.Range("N41").Value = '$/saco 25 kilos'
$ws.Range("O41").Value = 'Provincia del Elquí'
$ws.Range("P41").Value = 460

$ws.Range("D42").Value = 44473
$ws.Range("J42").Value = 500
$ws.Range("K42").Value = 8500
$ws.Range("L42").Value = 9000
$ws.Range("M42").Value = 8750
$ws.Range("N42").Value = '$/saco 25 kilos'
$ws.Range("O42").Value = 'Provincia del Elquí'
$ws.Range("P42").Value = 350

$ws.Range("D43").Value = 44874
$ws.Range("J43").Value = 500
$ws.Range("K43").Value = 6000
$ws.Range("L43").Value = 7000
$ws.Range("M43").Value = 6500
$ws.Range("N43").Value = '$/saco 25 kilos'
$ws.Range("O43").Value = 'Provincia del Elquí'
$ws.Range("P43").Value = 260
